$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC2's second step (row 20) currently holds the "liquidação" step content,
# and TC3's second step (row 28) currently holds the "atribuir/desatribuir"
# step content. Swap them so that:
#   - TC2 step 2 -> "atribuir/desatribuir" (Steps / Expected Results)
#   - TC3 step 2 -> "liquidação" (Steps / Expected Results)

$ws.Range("B20").Value = "Chefe Clica para atribuir/desatribuir o registro a si mesmo."
$ws.Range("D20").Value = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pela liquidação) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

$ws.Range("B28").Value = "Chefe Clica para realizar a liquidação."
$ws.Range("D28").Value = "SYSTEM Apresenta a tela de Registrar Liquidações"
